$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new "Wins" / "Losses" / "Ties" columns (AD:AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy header formatting (bold, centered, bordered) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record values for every data row (2-53)
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 72
    $ws.Cells.Item($r, 31).Value = 90
    $ws.Cells.Item($r, 32).Value = 0
}
